# Automatically re-generate list and index
# Review dates of the form "YYYY-MM-08" are shifted to "YYYY-MM-10".
# Values live in column B ("Review date") as plain text (not real dates).
# We force a Text number format just before writing so Excel's COM layer
# doesn't "helpfully" reinterpret the ISO-looking string as a date serial,
# then restore the cell's original (default/Normal) style so no visible
# formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -match '^(\d{4}-\d{2})-08$') {
        $newVal = "$($matches[1])-10"
        $cell.NumberFormat = "@"
        $cell.Value2 = $newVal
        $cell.Style = "Normal"
    }
}
